# Room-upload status sheet: add a small "Sheet2" that will be used to show
# whether an uploaded room record succeeded or failed, switch the workbook's
# base font from Sylfaen to Calibri, and refresh the active selection on
# both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Fix up the rich-text "ოთახის N" header on Sheet1 first, so the font
#     change is baked into the shared string *before* it gets copied onto
#     the new sheet below (keeps both sheets pointing at the same shared
#     string instead of forking a duplicate). -------------------------------
$headerChars = $ws1.Range("B1").Characters(7, 2)
$headerChars.Font.Name = "Calibri"

# --- Switch the workbook's default/"Normal" style font to Calibri. ---------
$normalStyle = $wb.Styles.Item("Normal")
$normalStyle.Font.Name = "Calibri"

# --- Refresh the selection left on Sheet1. ----------------------------------
[void]$ws1.Range("A14:C14").Select()

# --- Add the new sheet right after Sheet1. ----------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row - same three columns/headers as Sheet1.
$ws2.Range("A1").Value = "ტიპი"
$ws1.Range("B1").Copy($ws2.Range("B1"))
$ws2.Range("C1").Value = "სტუდ. ტევადობა"

# Single data row describing the newly uploaded room.
$ws2.Range("A2").Value = "მაგიდიანი"
$ws2.Range("B2").Value = "104ა"
$ws2.Range("C2").Value = 75

# Leave Sheet2 active with B3 selected, matching where the user would click
# next after reviewing the upload result.
[void]$ws2.Range("B3").Select()
$ws2.Activate()

Write-Host "done"
